# "use latest legacy import format of gbsl"
#
# The GBSL legacy import sheet used three generic placeholder columns
# ("Kategorie", "Kategorie2", "Kategorie3"). The latest legacy format
# renames these to the actual school-type short names used by GBSL:
#   Kategorie  -> GYM  (Gymnasium)
#   Kategorie2 -> FMS  (Fachmittelschule)
#   Kategorie3 -> WMS  (Wirtschaftsmittelschule)
#
# Renaming the header cells of the Excel table (Tabelle1) automatically
# renames the corresponding table columns (xl/tables/table1.xml) and
# updates the shared-string table accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "GYM"
$ws.Range("K1").Value = "FMS"
$ws.Range("L1").Value = "WMS"

# Restore the selection to the top of the re-labelled columns, as in the
# saved workbook (selection moved from R3 to G1).
[void]$ws.Range("G1").Select()
